$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on affected Price/Volume cells so numeric-looking strings are preserved as text
$textCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "D8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "E13", "D14", "E14", "D15", "E15", "E16", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "D26", "E26", "D38", "E38", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48", "E48", "D49", "E49", "D50", "E50", "D51", "E51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "307.29"
$ws.Range("E2").Value = "-2.57%"
$ws.Range("D3").Value = "41.03"
$ws.Range("E3").Value = "-2.17%"
$ws.Range("D4").Value = "5.050"
$ws.Range("E4").Value = "-2.66%"
$ws.Range("D5").Value = "0.07609"
$ws.Range("E5").Value = "-5.33%"
$ws.Range("D6").Value = "4.244"
$ws.Range("E6").Value = "-2.97%"
$ws.Range("D7").Value = "1.594"
$ws.Range("E7").Value = "-7.76%"
$ws.Range("D8").Value = "0.9058"
$ws.Range("D9").Value = "0.1005"
$ws.Range("E9").Value = "-10.46%"
$ws.Range("D10").Value = "0.1769"
$ws.Range("E10").Value = "-4.07%"
$ws.Range("D11").Value = "0.09070"
$ws.Range("E11").Value = "-1.55%"
$ws.Range("D12").Value = "0.04396"
$ws.Range("E12").Value = "-3.72%"
$ws.Range("E13").Value = "-0.04%"
$ws.Range("D14").Value = "0.001257"
$ws.Range("E14").Value = "-1.39%"
$ws.Range("D15").Value = "0.005875"
$ws.Range("E15").Value = "-1.40%"
$ws.Range("E16").Value = "0.39%"
$ws.Range("E17").Value = "-5.05%"
$ws.Range("D18").Value = "0.3297"
$ws.Range("E18").Value = "-2.58%"
$ws.Range("D19").Value = "6.840"
$ws.Range("E19").Value = "-7.00%"
$ws.Range("D20").Value = "0.1340"
$ws.Range("E20").Value = "-3.23%"
$ws.Range("E21").Value = "7.90%"
$ws.Range("D22").Value = "0.04153"
$ws.Range("E22").Value = "-0.80%"
$ws.Range("D23").Value = "0.001214"
$ws.Range("E23").Value = "-2.64%"
$ws.Range("D24").Value = "0.004005"
$ws.Range("E24").Value = "-5.57%"
$ws.Range("D25").Value = "0.0001301"
$ws.Range("E25").Value = "5.98%"
$ws.Range("D26").Value = "0.0003010"
$ws.Range("E26").Value = "0.61%"
$ws.Range("D38").Value = "0.02414"
$ws.Range("E38").Value = "-4.55%"
$ws.Range("E39").Value = "-4.57%"
$ws.Range("D40").Value = "0.007843"
$ws.Range("E40").Value = "-2.16%"
$ws.Range("D41").Value = "0.1305"
$ws.Range("E41").Value = "-6.03%"
$ws.Range("D42").Value = "0.007094"
$ws.Range("E42").Value = "-6.80%"
$ws.Range("D43").Value = "0.001950"
$ws.Range("E43").Value = "-5.83%"
$ws.Range("D44").Value = "0.008281"
$ws.Range("E44").Value = "-1.62%"
$ws.Range("D45").Value = "0.3056"
$ws.Range("E45").Value = "-2.62%"
$ws.Range("D46").Value = "0.00006385"
$ws.Range("E46").Value = "-5.58%"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").Value = "-0.33%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D48").Value = "0.005994"
$ws.Range("E48").Value = "76.36%"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D49").Value = "0.003003"
$ws.Range("E49").Value = "-27.07%"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").Value = "-0.33%"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").Value = "-0.33%"
